$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptos list refresh (prices in column D, 1h volume % in column E).
# Values that look like plain numbers are written with a leading apostrophe
# so Excel keeps them as literal text (matching the source feed's
# dotted/padded formatting) instead of coercing them into floating point.

$ws.Range("D2").Value = '36.057.38'
$ws.Range("E2").Value = '  -1.35%  '
$ws.Range("D3").Value = '1.920.62'
$ws.Range("E3").Value = '  -4.16%  '
$ws.Range("E4").Value = '  -0.13%  '
$ws.Range("D5").Value = '''239.09'
$ws.Range("E5").Value = '  -3.39%  '
$ws.Range("D6").Value = '''0.600'
$ws.Range("E6").Value = '  -4.91%  '
$ws.Range("E7").Value = '  -0.17%  '
$ws.Range("D8").Value = '''54.85'
$ws.Range("E8").Value = '  -11.14%  '
$ws.Range("D9").Value = '''0.358'
$ws.Range("E9").Value = '  -7.76%  '
$ws.Range("D10").Value = '''55.11'
$ws.Range("E10").Value = '  -3.13%  '
$ws.Range("D11").Value = '''0.0804'
$ws.Range("E11").Value = '  +3.60%  '
$ws.Range("E12").Value = '  -0.80%  '
$ws.Range("D13").Value = '2.207.47'
$ws.Range("E13").Value = '  -4.46%  '
$ws.Range("D14").Value = '''0.800'
$ws.Range("E14").Value = '  -8.79%  '
$ws.Range("D15").Value = '''20.52'
$ws.Range("E15").Value = '  -11.18%  '
$ws.Range("D16").Value = '''12.99'
$ws.Range("E16").Value = '  -8.25%  '
$ws.Range("D17").Value = '''5.12'
$ws.Range("E17").Value = '  -6.97%  '
$ws.Range("D18").Value = '1.922.81'
$ws.Range("E18").Value = '  -4.20%  '
$ws.Range("D19").Value = '35.946.95'
$ws.Range("E19").Value = '  -1.64%  '
$ws.Range("D20").Value = '''68.74'
$ws.Range("E20").Value = '  -4.26%  '
$ws.Range("D21").Value = '0.0₃0848'
$ws.Range("E21").Value = '  -2.65%  '
$ws.Range("D22").Value = '''224.66'
$ws.Range("E22").Value = '  -4.24%  '
$ws.Range("D23").Value = '''4.88'
$ws.Range("E23").Value = '  -7.55%  '
$ws.Range("E24").Value = '  +0.23%  '
$ws.Range("E25").Value = '  -4.84%  '
$ws.Range("D26").Value = '''2.23'
$ws.Range("E26").Value = '  -3.86%  '
$ws.Range("D27").Value = '''9.15'
$ws.Range("E27").Value = '  -5.95%  '
$ws.Range("D28").Value = '''161.71'
$ws.Range("E28").Value = '  +1.62%  '
$ws.Range("D29").Value = '''18.93'
$ws.Range("E29").Value = '  -5.48%  '
$ws.Range("E30").Value = '  -17.33%  '
$ws.Range("E31").Value = '  -3.36%  '
$ws.Range("D32").Value = '''1.11'
$ws.Range("E32").Value = '  -5.77%  '
$ws.Range("D33").Value = '''4.57'
$ws.Range("E33").Value = '  -7.72%  '
$ws.Range("D34").Value = '''0.0613'
$ws.Range("E34").Value = '  -0.54%  '
$ws.Range("D35").Value = '''4.18'
$ws.Range("E35").Value = '  -6.32%  '
$ws.Range("E36").Value = '  -0.13%  '
$ws.Range("D37").Value = '''1.79'
$ws.Range("E37").Value = '  -2.11%  '
$ws.Range("D38").Value = '''5.87'
$ws.Range("E38").Value = '  -9.92%  '
$ws.Range("D39").Value = '''2.10'
$ws.Range("E39").Value = '  -10.32%  '
$ws.Range("D40").Value = '''2.78'
$ws.Range("E40").Value = '  -11.74%  '
$ws.Range("D41").Value = '''0.0954'
$ws.Range("E41").Value = '  -3.98%  '
$ws.Range("D42").Value = '''2.86'
$ws.Range("E42").Value = '  -2.30%  '
$ws.Range("D43").Value = '''1.14'
$ws.Range("E43").Value = '  -8.84%  '
$ws.Range("D44").Value = '''0.0204'
$ws.Range("E44").Value = '  -4.25%  '
$ws.Range("D45").Value = '1.325.97'
$ws.Range("E45").Value = '  -2.03%  '
$ws.Range("D46").Value = '''15.25'
$ws.Range("E46").Value = '  -8.71%  '
$ws.Range("E47").Value = '  -9.83%  '
$ws.Range("D48").Value = '''86.02'
$ws.Range("E48").Value = '  -6.79%  '
$ws.Range("D49").Value = '''7.08'
$ws.Range("E49").Value = '  -6.94%  '
$ws.Range("E50").Value = '  -3.35%  '
$ws.Range("D51").Value = '''44.84'
$ws.Range("E51").Value = '  +0.64%  '
